$d = $word.ActiveDocument

# Locate the unique anchor text that ends the target paragraph (the list
# item about adding the professor's research content into the 志愿動機).
$range = $d.Content
$found = $range.Find.Execute("动机中加入学校教授的研究内容", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Move to the paragraph that contains the match, then collapse to the very
# end of that paragraph (just before the paragraph mark) so the new run is
# appended after the existing trailing whitespace runs.
$para = $range.Paragraphs(1)
$insertAt = $para.Range.Duplicate
$insertAt.MoveEnd(1, -1) | Out-Null
$insertAt.Collapse(0)

$insertAt.InsertAfter("待办")
$insertAt.Font.Name = "Cambria Math"
$insertAt.Font.NameAscii = "Cambria Math"
$insertAt.Font.NameFarEast = "Cambria Math"
$insertAt.HighlightColorIndex = 7
